# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps for the
# 78150835-813a-4f70-97fb-dc3c2db347f2 entries on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$zh = $wb.Worksheets.Item("zh-cn")
$zh.Cells.Item(3, 5).Value = "2016-03-25 07:00:01"
$zh.Cells.Item(3, 8).Value = "2016-03-25 07:00:49"
$zh.Cells.Item(5, 5).Value = "2016-03-25 07:00:01"
$zh.Cells.Item(5, 8).Value = "2016-03-25 07:00:49"

$de = $wb.Worksheets.Item("de-de")
$de.Cells.Item(3, 5).Value = "2016-03-25 07:00:10"
$de.Cells.Item(3, 8).Value = "2016-03-25 07:00:57"
$de.Cells.Item(5, 5).Value = "2016-03-25 07:00:10"
$de.Cells.Item(5, 8).Value = "2016-03-25 07:00:57"
